$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows belonging to "BRIAN RAFAEL DELGADO SIERRA" (row 20)
# and the duplicated "MAURA ALEJANDRA VARON SIERRA" rows (rows 21-22).
# Deleting shifts the trailing rows (old row 23 -> 20, old rows 28-29 -> 25-26)
# up automatically, carrying their formatting/styles along.
$ws.Rows("20:22").Delete()

# Update the account summary figures.
$ws.Range("E11").Value = 217862
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 4

# Row 16: CARLOS ALBERTO LARA MARTINEZ -> MAURA ALEJANDRA VARON SIERRA (period 2207)
$ws.Range("C16").Value = "1047485077"
$ws.Range("D16").Value = "MAURA ALEJANDRA VARON SIERRA"
$ws.Range("E16").Value = "2207"
$ws.Range("F16").Value = 64000
$ws.Range("G16").Value = 2300000

# Row 17: CARLOS ALBERTO LARA MARTINEZ -> MAURA ALEJANDRA VARON SIERRA (period 2211)
$ws.Range("C17").Value = "1047485077"
$ws.Range("D17").Value = "MAURA ALEJANDRA VARON SIERRA"
$ws.Range("E17").Value = "2211"
$ws.Range("F17").Value = 64000
$ws.Range("G17").Value = 2300000

# Row 18: YEINER DE JESUS PEREZ RIVERA -> NICHOLLE LOPEZ GALINDO (period 2401)
$ws.Range("C18").Value = "1002198300"
$ws.Range("D18").Value = "NICHOLLE LOPEZ GALINDO"
$ws.Range("E18").Value = "2401"
$ws.Range("F18").Value = 1733
$ws.Range("G18").Value = 1562473

# Row 19: NICHOLLE LOPEZ GALINDO -> YEINER DE JESUS PEREZ RIVERA (period 2503)
$ws.Range("C19").Value = "1143384255"
$ws.Range("D19").Value = "YEINER DE JESUS PEREZ RIVERA"
$ws.Range("E19").Value = "2503"
$ws.Range("F19").Value = 85892
$ws.Range("G19").Value = 2147300

# Row 20 (formerly row 23, RICHARD VILLAMORO MORELOS): minor value correction.
$ws.Range("G20").Value = 1677580
